$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text columns (Coin name / Link URL) - direct assignment ---
$ws.Range('B7').Value = 'GateToken'
$ws.Range('C7').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('B8').Value = 'FTXToken'
$ws.Range('C8').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('B9').Value = 'MXToken'
$ws.Range('C9').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('B10').Value = 'WazirX'
$ws.Range('C10').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('B11').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C11').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('B12').Value = 'MandalaExchangeToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('B13').Value = 'BitrueCoin'
$ws.Range('C13').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('B14').Value = 'BitMartToken'
$ws.Range('C14').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('B15').Value = 'BitForexToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('B16').Value = 'One'
$ws.Range('C16').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('B17').Value = 'TigerCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('B18').Value = 'LEO'
$ws.Range('C18').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'

# --- Numeric-looking text columns (Price / Volume%) ---
# Force text storage (matches source inlineStr cells) by temporarily
# formatting as Text, writing the literal string, then clearing the
# format again so the cell keeps the workbook default style.
$textCells = @(
    'D2', 'E2', 'D3', 'E3', 'D4', 'E4', 'D5', 'E5', 'D6', 'E6', 'D7', 'E7', 'D8', 'E8', 'D9', 'E9', 'D10', 'E10', 'D11', 'E11', 'D12', 'E12', 'D13', 'E13', 'D14', 'E14', 'D15', 'E15', 'D16', 'E16', 'D17', 'E17', 'D18', 'E18', 'E19', 'E20', 'D21', 'E21', 'D22', 'E22', 'D23', 'E23', 'D24', 'E24', 'D25', 'E25', 'D26', 'E26', 'D27', 'E27', 'D28', 'E28', 'D40', 'E40', 'D41', 'E41', 'D42', 'E42', 'D43', 'E43', 'D44', 'E44', 'D45', 'E45', 'E46', 'D47', 'E47'
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '287.18'
$ws.Range('E2').Value = '1.84%'
$ws.Range('D3').Value = '29.40'
$ws.Range('E3').Value = '4.36%'
$ws.Range('D4').Value = '5.104'
$ws.Range('E4').Value = '1.36%'
$ws.Range('D5').Value = '0.06851'
$ws.Range('E5').Value = '5.60%'
$ws.Range('D6').Value = '7.367'
$ws.Range('E6').Value = '1.69%'
$ws.Range('D7').Value = '3.513'
$ws.Range('E7').Value = '4.25%'
$ws.Range('D8').Value = '1.383'
$ws.Range('E8').Value = '-0.13%'
$ws.Range('D9').Value = '0.9031'
$ws.Range('E9').Value = '-2.78%'
$ws.Range('D10').Value = '0.1592'
$ws.Range('E10').Value = '2.67%'
$ws.Range('D11').Value = '0.07041'
$ws.Range('E11').Value = '13.39%'
$ws.Range('D12').Value = '0.07617'
$ws.Range('E12').Value = '1.16%'
$ws.Range('D13').Value = '0.02920'
$ws.Range('E13').Value = '1.64%'
$ws.Range('D14').Value = '0.08985'
$ws.Range('E14').Value = '0.10%'
$ws.Range('D15').Value = '0.001585'
$ws.Range('E15').Value = '-0.36%'
$ws.Range('D16').Value = '0.0006506'
$ws.Range('E16').Value = '2.09%'
$ws.Range('D17').Value = '0.006390'
$ws.Range('E17').Value = '5.22%'
$ws.Range('D18').Value = '3.470'
$ws.Range('E18').Value = '0.83%'
$ws.Range('E19').Value = '0.06%'
$ws.Range('E20').Value = '0.83%'
$ws.Range('D21').Value = '0.1323'
$ws.Range('E21').Value = '1.48%'
$ws.Range('D22').Value = '4.014'
$ws.Range('E22').Value = '-1.28%'
$ws.Range('D23').Value = '0.1553'
$ws.Range('E23').Value = '0.41%'
$ws.Range('D24').Value = '0.04499'
$ws.Range('E24').Value = '1.54%'
$ws.Range('D25').Value = '0.001204'
$ws.Range('E25').Value = '1.78%'
$ws.Range('D26').Value = '0.004380'
$ws.Range('E26').Value = '-0.25%'
$ws.Range('D27').Value = '0.0001165'
$ws.Range('E27').Value = '-6.71%'
$ws.Range('D28').Value = '0.0001611'
$ws.Range('E28').Value = '-0.53%'
$ws.Range('D40').Value = '0.04290'
$ws.Range('E40').Value = '3.34%'
$ws.Range('D41').Value = '0.006786'
$ws.Range('E41').Value = '2.11%'
$ws.Range('D42').Value = '0.1245'
$ws.Range('E42').Value = '2.01%'
$ws.Range('D43').Value = '0.002181'
$ws.Range('E43').Value = '8.07%'
$ws.Range('D44').Value = '0.01160'
$ws.Range('E44').Value = '-4.18%'
$ws.Range('D45').Value = '0.00005727'
$ws.Range('E45').Value = '2.46%'
$ws.Range('E46').Value = '-1.86%'
$ws.Range('D47').Value = '0.01301'
$ws.Range('E47').Value = '0.00%'

foreach ($addr in $textCells) {
    $ws.Range($addr).ClearFormats()
}
